$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at position 31, shifting existing rows 31-78 down to 32-79.
$ws.Rows("31:31").Insert()

# Populate the new row's September Details / Date cells.
$ws.Range("R31").Value = "bal axisbank"
$ws.Range("S31").Value = "2024-09-05 16:52:25"
